# Refresh the crypto price/volume snapshot (Price = column D, Volume(1h) =
# column E) for rows 2-51 with values from the latest scrape.
#
# Column D cells are plain text (e.g. "29.585.51", "0.9994"), not real
# numbers -- the source site uses "." as a thousands separator in some
# rows, so the whole column is stored as text. Setting .Value directly
# with a string that *looks* numeric (e.g. "1.000", "0.9994") would make
# Excel silently reinterpret it as a Number and drop the trailing zeros,
# so for those we prefix with a literal apostrophe first (the same
# quote-prefix text-entry trick Excel itself uses) to force a text cell.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "29.585.51"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "1.856.81"
$ws.Range("E3").Value = "  +1.33%  "
Set-TextValue $ws.Range("D4") "0.9994"
Set-TextValue $ws.Range("D5") "244.78"
$ws.Range("E5").Value = "  -0.20%  "
Set-TextValue $ws.Range("D6") "0.6939"
$ws.Range("E6").Value = "  +0.52%  "
Set-TextValue $ws.Range("D7") "1.000"
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue $ws.Range("D8") "0.07689"
$ws.Range("E8").Value = "  +0.48%  "
Set-TextValue $ws.Range("D9") "0.3061"
$ws.Range("E9").Value = "  +0.25%  "
Set-TextValue $ws.Range("D10") "23.67"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "1.858.92"
$ws.Range("E12").Value = "  +1.41%  "
Set-TextValue $ws.Range("D13") "5.140"
$ws.Range("E13").Value = "  +1.51%  "
Set-TextValue $ws.Range("D14") "91.03"
$ws.Range("E14").Value = "  +0.76%  "
Set-TextValue $ws.Range("D15") "0.6910"
$ws.Range("E15").Value = "  +1.86%  "
Set-TextValue $ws.Range("D16") "6.523"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "29.471.42"
$ws.Range("E17").Value = "  +1.92%  "
Set-TextValue $ws.Range("D18") "0.000008284"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "2.106.31"
$ws.Range("E19").Value = "  +1.22%  "
Set-TextValue $ws.Range("D20") "238.42"
$ws.Range("E20").Value = "  -1.71%  "
Set-TextValue $ws.Range("D21") "12.73"
$ws.Range("E21").Value = "  +0.45%  "
Set-TextValue $ws.Range("D22") "0.9997"
$ws.Range("E22").Value = "  +0.05%  "
Set-TextValue $ws.Range("D23") "7.607"
$ws.Range("E23").Value = "  +2.22%  "
Set-TextValue $ws.Range("D24") "1.000"
$ws.Range("E24").Value = "  +0.07%  "
Set-TextValue $ws.Range("D25") "0.1493"
$ws.Range("E25").Value = "  +1.83%  "
Set-TextValue $ws.Range("D26") "8.898"
$ws.Range("E26").Value = "  +1.16%  "
Set-TextValue $ws.Range("D27") "159.38"
$ws.Range("E27").Value = "  -0.96%  "
Set-TextValue $ws.Range("D28") "18.23"
$ws.Range("E28").Value = "  +0.22%  "
Set-TextValue $ws.Range("D29") "1.540"
$ws.Range("E29").Value = "  -1.23%  "
Set-TextValue $ws.Range("D30") "4.248"
$ws.Range("E30").Value = "  +0.77%  "
Set-TextValue $ws.Range("D31") "4.165"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +2.52%  "
Set-TextValue $ws.Range("D33") "0.05155"
$ws.Range("E33").Value = "  +0.80%  "
Set-TextValue $ws.Range("D34") "0.7685"
$ws.Range("E34").Value = "  +1.73%  "
Set-TextValue $ws.Range("D35") "1.888"
$ws.Range("E35").Value = "  +2.47%  "
Set-TextValue $ws.Range("D36") "1.150"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "1.332.02"
$ws.Range("E38").Value = "  +8.51%  "
Set-TextValue $ws.Range("D39") "0.01870"
$ws.Range("E39").Value = "  +1.57%  "
Set-TextValue $ws.Range("D40") "2.725"
$ws.Range("E40").Value = "  +1.34%  "
Set-TextValue $ws.Range("D41") "0.9705"
$ws.Range("E41").Value = "  +4.69%  "
Set-TextValue $ws.Range("D42") "106.57"
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("E43").Value = "  +1.92%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").Value = "2.003.89"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("E47").Value = "  +1.72%  "
Set-TextValue $ws.Range("D48") "0.5217"
$ws.Range("E48").Value = "  +0.84%  "
Set-TextValue $ws.Range("D49") "1.775"
$ws.Range("E49").Value = "  +2.09%  "
Set-TextValue $ws.Range("D50") "62.81"
$ws.Range("E50").Value = "  -2.03%  "
Set-TextValue $ws.Range("D51") "6.949"
$ws.Range("E51").Value = "  +0.57%  "
